$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stream_me")

$ws.Range("C48").Value = "slow"

$ws.Range("C49").Value = "yellow"
$ws.Range("D49").Value = 2
$ws.Range("F49").Value = 0.1564945226917058

$ws.Range("C50").Value = "dumb"
$ws.Range("D50").Value = 1
$ws.Range("F50").Value = 0.0782472613458529

$ws.Range("C51").Value = "blessed"

$ws.Range("C52").Value = "classic"

$ws.Range("C53").Value = "laffy"

$ws.Range("C54").Value = "reassure"

$ws.Range("C77").Value = "classic"
$ws.Range("D77").Value = 1
$ws.Range("F77").Value = 0.07215007215007214

$ws.Range("C78").Value = "blessed"

$ws.Range("D132").Value = 1
$ws.Range("F132").Value = 0.08064516129032258

$ws.Range("C133").Value = "stanky"
$ws.Range("D133").Value = 1
$ws.Range("F133").Value = 0.08064516129032258

$ws.Range("C134").Value = "dumb"

$ws.Range("C135").Value = "laffy"

$ws.Range("C136").Value = "scared"

$ws.Range("C160").Value = "classic"
$ws.Range("D160").Value = 1
$ws.Range("F160").Value = 0.08873114463176575

$ws.Range("C161").Value = "blessed"

$ws.Range("C186").Value = "dumb"
$ws.Range("D186").Value = 1
$ws.Range("F186").Value = 0.0999000999000999

$ws.Range("C187").Value = "reassure"
$ws.Range("D187").Value = 1
$ws.Range("F187").Value = 0.0999000999000999

$ws.Range("C188").Value = "blessed"

$ws.Range("C189").Value = "classic"

$ws.Range("C190").Value = "laffy"

$ws.Range("C212").Value = "dumb"
$ws.Range("D212").Value = 1
$ws.Range("F212").Value = 0.08710801393728224

$ws.Range("C213").Value = "blessed"

$ws.Range("C214").Value = "classic"

$ws.Range("C240").Value = "reassure"
$ws.Range("D240").Value = 1
$ws.Range("F240").Value = 0.08665511265164645

$ws.Range("C241").Value = "blessed"

$ws.Range("C242").Value = "classic"

$ws.Range("C243").Value = "laffy"

$ws.Range("C292").Value = "classic"
$ws.Range("D292").Value = 1
$ws.Range("F292").Value = 0.06901311249137336

$ws.Range("C293").Value = "dumb"
$ws.Range("D293").Value = 1
$ws.Range("F293").Value = 0.06901311249137336

$ws.Range("C294").Value = "blessed"

$ws.Range("C320").Value = "slow"
$ws.Range("D320").Value = 2
$ws.Range("F320").Value = 0.1573564122738002

$ws.Range("C321").Value = "yellow"
$ws.Range("D321").Value = 1
$ws.Range("F321").Value = 0.07867820613690008

$ws.Range("C322").Value = "blessed"

$ws.Range("C323").Value = "laffy"

$ws.Range("C324").Value = "reassure"

$ws.Range("C325").Value = "stanky"

$ws.Range("C400").Value = "slow"

$ws.Range("C401").Value = "yellow"
$ws.Range("D401").Value = 2
$ws.Range("F401").Value = 0.1616814874696847

$ws.Range("C402").Value = "blessed"

$ws.Range("C403").Value = "classic"

$ws.Range("C404").Value = "laffy"

$ws.Range("C405").Value = "scared"

$ws.Range("C428").Value = "scared"
$ws.Range("D428").Value = 1
$ws.Range("F428").Value = 0.08077544426494346

$ws.Range("C429").Value = "stanky"
$ws.Range("D429").Value = 1
$ws.Range("F429").Value = 0.08077544426494346

$ws.Range("C430").Value = "blessed"

$ws.Range("C431").Value = "classic"

$ws.Range("C432").Value = "dumb"

$ws.Range("C433").Value = "laffy"
